$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 2561.125
$ws.Range("I5").Value = 2918.4285
$ws.Range("K5").Value = 2918.4285
$ws.Range("M5").Value = -2803.4285
$ws.Range("H33").Value = 265.56522
$ws.Range("I33").Value = 196.73685
$ws.Range("K33").Value = 196.73685
$ws.Range("M33").Value = 32.26315
$ws.Range("H38").Value = 398.5
$ws.Range("J38").Value = 1000
$ws.Range("L38").Value = 3000
$ws.Range("N38").Value = -3744
$ws.Range("H107").Value = 594.2105
$ws.Range("I107").Value = 507.77777
$ws.Range("J107").Value = 672
$ws.Range("K107").Value = 507.77777
$ws.Range("L107").Value = 672
$ws.Range("M107").Value = 1412.22223
$ws.Range("N107").Value = -4512
$ws.Range("H111").Value = 884.1429000000001
$ws.Range("I111").Value = 832
$ws.Range("J111").Value = 953.6667
$ws.Range("K111").Value = 2496
$ws.Range("L111").Value = 2861.0001
$ws.Range("M111").Value = 571
$ws.Range("N111").Value = -8995.000100000001
$ws.Range("H115").Value = 137.5
$ws.Range("I115").Value = 137.5
$ws.Range("K115").Value = 412.5
$ws.Range("M115").Value = 1154.5
$ws.Range("H123").Value = 151999.5
$ws.Range("J123").Value = 151999.5
$ws.Range("L123").Value = 151999.5
$ws.Range("N123").Value = -161799.5
$ws.Range("H125").Value = 171262.67
$ws.Range("I125").Value = 3919.5
$ws.Range("K125").Value = 35275.5
$ws.Range("M125").Value = -32815.5
$ws.Range("H132").Value = 3999.875
$ws.Range("I132").Value = 4457
$ws.Range("K132").Value = 13371
$ws.Range("M132").Value = -10841
$ws.Range("H138").Value = 2362
$ws.Range("J138").Value = 2934.375
$ws.Range("L138").Value = 8803.125
$ws.Range("N138").Value = -19083.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 255.66667
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 255.66667
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 255.66667
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -487.66667
$ws.Range("H5").Value = 475
$ws.Range("I5").Value = 475
$ws.Range("K5").Value = 475
$ws.Range("M5").Value = -363
$ws.Range("H6").Value = 4000
$ws.Range("J6").Value = 4000
$ws.Range("L6").Value = 4000
$ws.Range("N6").Value = -4346
$ws.Range("H46").Value = 3625.5
$ws.Range("J46").Value = 3625.5
$ws.Range("L46").Value = 3625.5
$ws.Range("N46").Value = -4263.5
$ws.Range("H55").Value = 27500
$ws.Range("J55").Value = 27500
$ws.Range("L55").Value = 27500
$ws.Range("N55").Value = -28130
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984
$ws.Range("H97").Value = 789.8
$ws.Range("I97").Value = 789.8
$ws.Range("K97").Value = 789.8
$ws.Range("M97").Value = -293.8
$ws.Range("H110").Value = 3649
$ws.Range("I110").Value = 1850
$ws.Range("K110").Value = 1850
$ws.Range("M110").Value = 195
$ws.Range("H122").Value = 2201.0908
$ws.Range("I122").Value = 1289.5
$ws.Range("K122").Value = 3868.5
$ws.Range("M122").Value = -1418.5
$ws.Range("H132").Value = 3012
$ws.Range("I132").Value = 3012
$ws.Range("K132").Value = 9036
$ws.Range("M132").Value = -6506

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 475
$ws.Range("I4").Value = 475
$ws.Range("K4").Value = 475
$ws.Range("M4").Value = -360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 21000
$ws.Range("J70").Value = 21000
$ws.Range("L70").Value = 21000
$ws.Range("N70").Value = -21630
$ws.Range("H73").Value = 21000
$ws.Range("J73").Value = 21000
$ws.Range("L73").Value = 21000
$ws.Range("N73").Value = -23184
$ws.Range("H105").Value = 1234.2858
$ws.Range("I105").Value = 994
$ws.Range("J105").Value = 1835
$ws.Range("K105").Value = 994
$ws.Range("L105").Value = 1835
$ws.Range("M105").Value = 753
$ws.Range("N105").Value = -5329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 177.63637
$ws.Range("J2").Value = 34.666668
$ws.Range("L2").Value = 208.000008
$ws.Range("N2").Value = -434.000008
$ws.Range("H17").Value = 825.8
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 1029.75
$ws.Range("K17").Value = 30
$ws.Range("L17").Value = 3089.25
$ws.Range("M17").Value = 139
$ws.Range("N17").Value = -3427.25
$ws.Range("H34").Value = 1482.75
$ws.Range("I34").Value = 500
$ws.Range("J34").Value = 1572.091
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 4716.272999999999
$ws.Range("M34").Value = -1416
$ws.Range("N34").Value = -4884.272999999999
$ws.Range("H39").Value = 5749.6665
$ws.Range("J39").Value = 5749.6665
$ws.Range("L39").Value = 17248.9995
$ws.Range("N39").Value = -17836.9995
$ws.Range("H55").Value = 2072.2727
$ws.Range("I55").Value = 595
$ws.Range("J55").Value = 2220
$ws.Range("K55").Value = 1785
$ws.Range("L55").Value = 6660
$ws.Range("M55").Value = -1608
$ws.Range("N55").Value = -7014
$ws.Range("H131").Value = 1995
$ws.Range("I131").Value = 1995
$ws.Range("K131").Value = 5985
$ws.Range("M131").Value = -945
$ws.Range("H136").Value = 3025
$ws.Range("I136").Value = 3025
$ws.Range("K136").Value = 9075
$ws.Range("M136").Value = -3975

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 25000
$ws.Range("J47").Value = 25000
$ws.Range("L47").Value = 25000
$ws.Range("N47").Value = -26136
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()
$ws.Range("H107").Value = 154.66667
$ws.Range("I107").Value = 160.8
$ws.Range("J107").Value = 124
$ws.Range("K107").Value = 160.8
$ws.Range("L107").Value = 124
$ws.Range("M107").Value = 1759.2
$ws.Range("N107").Value = -3964
$ws.Range("H122").Value = 1399.7142
$ws.Range("I122").Value = 1450
$ws.Range("J122").Value = 1274
$ws.Range("K122").Value = 4350
$ws.Range("L122").Value = 3822
$ws.Range("M122").Value = -1900
$ws.Range("N122").Value = -8722
$ws.Range("H126").Value = 16777.75
$ws.Range("I126").Value = 12370.333
$ws.Range("J126").Value = 30000
$ws.Range("K126").Value = 37110.999
$ws.Range("L126").Value = 90000
$ws.Range("M126").Value = -34640.999
$ws.Range("N126").Value = -94940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 20195.154
$ws.Range("I7").Value = 19317.75
$ws.Range("J7").Value = 21599
$ws.Range("K7").Value = 19317.75
$ws.Range("L7").Value = 21599
$ws.Range("M7").Value = -19205.75
$ws.Range("N7").Value = -21823
$ws.Range("H40").Value = 7505.25
$ws.Range("I40").Value = 7505.25
$ws.Range("K40").Value = 7505.25
$ws.Range("M40").Value = -7369.25
$ws.Range("H46").Value = 3882.353
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 4875
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 4875
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -5251
$ws.Range("H55").Value = 212.27272
$ws.Range("I55").Value = 170.14285
$ws.Range("K55").Value = 170.14285
$ws.Range("M55").Value = 2.85714999999999
$ws.Range("H82").Value = 210
$ws.Range("I82").Value = 210
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 210
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = 151
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 210
$ws.Range("I85").Value = 210
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 210
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 1038
$ws.Range("N85").ClearContents()
$ws.Range("H126").Value = 20195.154
$ws.Range("I126").Value = 19317.75
$ws.Range("J126").Value = 21599
$ws.Range("K126").Value = 57953.25
$ws.Range("L126").Value = 64797
$ws.Range("M126").Value = -55483.25
$ws.Range("N126").Value = -69737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 8527.5
$ws.Range("I32").Value = 7026
$ws.Range("K32").Value = 7026
$ws.Range("M32").Value = -6709
$ws.Range("H54").Value = 16250
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 50000
$ws.Range("N65").Value = -56240
$ws.Range("H70").Value = 32857.145
$ws.Range("J70").Value = 32857.145
$ws.Range("L70").Value = 32857.145
$ws.Range("N70").Value = -33487.145
$ws.Range("H73").Value = 32857.145
$ws.Range("J73").Value = 32857.145
$ws.Range("L73").Value = 32857.145
$ws.Range("N73").Value = -35041.145
$ws.Range("H75").Value = 35000
$ws.Range("J75").Value = 35000
$ws.Range("L75").Value = 35000
$ws.Range("N75").Value = -36872
$ws.Range("H78").Value = 35000
$ws.Range("J78").Value = 35000
$ws.Range("L78").Value = 105000
$ws.Range("N78").Value = -114360
$ws.Range("H132").Value = 2603.6667
$ws.Range("I132").Value = 2494.4
$ws.Range("J132").Value = 3150
$ws.Range("K132").Value = 7483.200000000001
$ws.Range("L132").Value = 9450
$ws.Range("M132").Value = -4953.200000000001
$ws.Range("N132").Value = -14510

